$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change) per upstream refresh.
$ws.Range("D2").Value = "29.437.03"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.850.25"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "240.41"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.07675"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("D9").Value = "0.2922"
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").Value = "24.91"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("D11").Value = "0.07754"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "1.843.64"
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("D13").Value = "5.037"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6820"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("E15").Value = "  +2.14%  "
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "6.221"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "29.446.81"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "228.76"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.450"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "157.92"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1380"
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("D26").Value = "8.421"
$ws.Range("E26").Value = "  +0.93%  "
$ws.Range("D28").Value = "1.373"
$ws.Range("E28").Value = "  +5.54%  "
$ws.Range("D29").Value = "1.459"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "0.05628"
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("D31").Value = "4.128"
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("D33").Value = "1.845"
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("D34").Value = "1.164"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "0.7072"
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("D37").Value = "1.225.66"
$ws.Range("E37").Value = "  -1.92%  "
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("D39").Value = "2.757"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").Value = "6.448"
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").Value = "0.9028"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "2.011.29"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("D44").Value = "101.99"
$ws.Range("D45").Value = "66.19"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").Value = "7.209"
$ws.Range("E46").Value = "  +1.66%  "
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("D48").Value = "0.4021"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").Value = "9.039"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("E50").Value = "  +2.91%  "
$ws.Range("E51").Value = "  -0.30%  "
